$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1: A1 holds the text label "HK_G_acc_G" (shared string table now has a duplicate entry inserted
# at the front in the source diff; the cell's displayed/string value is unchanged).
$ws.Range("A1").Value = "HK_G_acc_G"

# Rows 2-49: updated numeric threshold values
$ws.Cells.Item(2, 1).Value = 50.070323488045013
$ws.Cells.Item(3, 1).Value = 49.929676511954995
$ws.Cells.Item(4, 1).Value = 50.210970464135016
$ws.Cells.Item(5, 1).Value = 49.507735583684955
$ws.Cells.Item(6, 1).Value = 49.507735583684955
$ws.Cells.Item(7, 1).Value = 49.367088607594937
$ws.Cells.Item(8, 1).Value = 51.61744022503516
$ws.Cells.Item(9, 1).Value = 51.336146272855132
$ws.Cells.Item(10, 1).Value = 51.61744022503516
$ws.Cells.Item(11, 1).Value = 51.898734177215189
$ws.Cells.Item(12, 1).Value = 50.351617440225041
$ws.Cells.Item(13, 1).Value = 50.914205344585092
$ws.Cells.Item(14, 1).Value = 51.758087201125178
$ws.Cells.Item(15, 1).Value = 51.61744022503516
$ws.Cells.Item(16, 1).Value = 52.180028129395218
$ws.Cells.Item(17, 1).Value = 50.492264416315045
$ws.Cells.Item(18, 1).Value = 50.773558368495074
$ws.Cells.Item(19, 1).Value = 49.929676511954995
$ws.Cells.Item(20, 1).Value = 50.632911392405063
$ws.Cells.Item(21, 1).Value = 51.054852320675103
$ws.Cells.Item(22, 1).Value = 50.914205344585092
$ws.Cells.Item(23, 1).Value = 48.804500703234879
$ws.Cells.Item(24, 1).Value = 48.38255977496484
$ws.Cells.Item(25, 1).Value = 48.52320675105485
$ws.Cells.Item(26, 1).Value = 50.632911392405063
$ws.Cells.Item(27, 1).Value = 50.914205344585092
$ws.Cells.Item(28, 1).Value = 51.195499296765121
$ws.Cells.Item(29, 1).Value = 52.180028129395218
$ws.Cells.Item(30, 1).Value = 51.47679324894515
$ws.Cells.Item(31, 1).Value = 51.758087201125178
$ws.Cells.Item(32, 1).Value = 48.945147679324897
$ws.Cells.Item(33, 1).Value = 49.367088607594937
$ws.Cells.Item(34, 1).Value = 49.789029535864984
$ws.Cells.Item(35, 1).Value = 50.632911392405063
$ws.Cells.Item(36, 1).Value = 50.070323488045013
$ws.Cells.Item(37, 1).Value = 53.586497890295362
$ws.Cells.Item(38, 1).Value = 48.945147679324897
$ws.Cells.Item(39, 1).Value = 49.929676511954995
$ws.Cells.Item(40, 1).Value = 49.929676511954995
$ws.Cells.Item(41, 1).Value = 51.47679324894515
$ws.Cells.Item(42, 1).Value = 51.61744022503516
$ws.Cells.Item(43, 1).Value = 51.758087201125178
$ws.Cells.Item(44, 1).Value = 51.61744022503516
$ws.Cells.Item(45, 1).Value = 50.351617440225041
$ws.Cells.Item(46, 1).Value = 50.351617440225041
$ws.Cells.Item(47, 1).Value = 49.226441631504926
$ws.Cells.Item(48, 1).Value = 51.898734177215189
$ws.Cells.Item(49, 1).Value = 50.632911392405063
